# Rebuild the "Лист1" project-tracking table with the new header/content
# and two additional student rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("A1").Value = "ФИО Преподавателя"
$ws.Range("B1").Value = "Название проекта"
$ws.Range("C1").Value = "Фамилия студента"
$ws.Range("D1").Value = "Имя студента"
$ws.Range("E1").Value = "Группа студента"
$ws.Range("F1").Value = "Статус"
$ws.Range("G1").Value = "Оценка"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = "Винаева Дарья"
$ws.Range("B2").Value = "гойда"
$ws.Range("C2").Value = "Нетеса"
$ws.Range("D2").Value = "Роман"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3433"
$ws.Range("F2").Value = "активен"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5"

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "Винаева Дарья"
$ws.Range("B3").Value = "гойда 2"
$ws.Range("C3").Value = "Нетеса"
$ws.Range("D3").Value = "Роман"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3433"
$ws.Range("F3").Value = "в разработке"
$ws.Range("G3").Value = "не оценено"

# --- Row 4 (new) ----------------------------------------------------------
$ws.Range("A4").Value = "Винаева Дарья"
$ws.Range("B4").Value = "гойда 2"
$ws.Range("C4").Value = "Макарова"
$ws.Range("D4").Value = "Ольга"
$ws.Range("E4").Value = 4317
$ws.Range("F4").Value = "в разработке"
$ws.Range("G4").Value = "не оценено"

# --- Row 5 (new) ----------------------------------------------------------
$ws.Range("A5").Value = "Винаева Дарья"
$ws.Range("B5").Value = "проверка"
$ws.Range("C5").Value = "Иванов"
$ws.Range("D5").Value = "Иван"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4319"
$ws.Range("F5").Value = "активен"
$ws.Range("G5").Value = "не оценено"

# --- Column A width ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13

# --- Selection / active cell ---------------------------------------------
$ws.Range("A5:G5").Select()
